# Apply cryptocurrency price/volume updates to the cryptos worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.987.05"
$ws.Range("E2").Value = "  -3.01%  "

$ws.Range("D3").Value = "3.332.42"
$ws.Range("E3").Value = "  -4.93%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'551.53"
$ws.Range("E5").Value = "  -4.49%  "

$ws.Range("D6").Value = "'172.48"
$ws.Range("E6").Value = "  -3.24%  "

$ws.Range("D7").Value = "'0.612"
$ws.Range("E7").Value = "  -3.55%  "

$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").Value = "3.326.68"
$ws.Range("E9").Value = "  -4.88%  "

$ws.Range("D10").Value = "'0.622"
$ws.Range("E10").Value = "  -1.92%  "

$ws.Range("E11").Value = "  +2.76%  "

$ws.Range("D12").Value = "'53.00"
$ws.Range("E12").Value = "  -4.37%  "

$ws.Range("D13").Value = "'0.0000271"
$ws.Range("E13").Value = "  -0.58%  "

$ws.Range("D14").Value = "'8.99"
$ws.Range("E14").Value = "  -2.62%  "

$ws.Range("D15").Value = "3.864.78"
$ws.Range("E15").Value = "  -5.03%  "

$ws.Range("D16").Value = "'18.21"
$ws.Range("E16").Value = "  -0.85%  "

$ws.Range("E17").Value = "  -3.52%  "

$ws.Range("D18").Value = "3.332.82"
$ws.Range("E18").Value = "  -4.92%  "

$ws.Range("D19").Value = "'11.76"
$ws.Range("E19").Value = "  -2.55%  "

$ws.Range("D20").Value = "63.879.90"
$ws.Range("E20").Value = "  -3.08%  "

$ws.Range("D21").Value = "'0.970"
$ws.Range("E21").Value = "  -3.44%  "

$ws.Range("D22").Value = "'427.05"
$ws.Range("E22").Value = "  +3.06%  "

$ws.Range("D23").Value = "'4.65"
$ws.Range("E23").Value = "  +7.65%  "

$ws.Range("D24").Value = "'4.08"
$ws.Range("E24").Value = "  -3.37%  "

$ws.Range("D25").Value = "'13.54"
$ws.Range("E25").Value = "  +1.96%  "

$ws.Range("D26").Value = "'84.11"
$ws.Range("E26").Value = "  -2.04%  "

$ws.Range("D27").Value = "'10.63"
$ws.Range("E27").Value = "  -3.45%  "

$ws.Range("E28").Value = "  -1.25%  "

$ws.Range("D29").Value = "'8.56"
$ws.Range("E29").Value = "  -6.03%  "

$ws.Range("D30").Value = "'29.62"
$ws.Range("E30").Value = "  -2.67%  "

$ws.Range("E31").Value = "  +1.74%  "

$ws.Range("D32").Value = "'594.15"
$ws.Range("E32").Value = "  -5.27%  "

$ws.Range("D33").Value = "'11.40"
$ws.Range("E33").Value = "  -2.38%  "

$ws.Range("E34").Value = "  -3.22%  "

$ws.Range("D35").Value = "'58.12"
$ws.Range("E35").Value = "  -2.47%  "

$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("E37").Value = "  -8.77%  "

$ws.Range("E38").Value = "  +0.45%  "

$ws.Range("D39").Value = "'35.26"
$ws.Range("E39").Value = "  -5.18%  "

$ws.Range("D40").Value = "0.0₃0746"
$ws.Range("E40").Value = "  -6.43%  "

$ws.Range("E41").Value = "  -4.52%  "

$ws.Range("D42").Value = "3.101.38"
$ws.Range("E42").Value = "  -5.29%  "

$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.04%  "

$ws.Range("E44").Value = "  -5.01%  "

$ws.Range("D45").Value = "'0.0404"
$ws.Range("E45").Value = "  -3.13%  "

$ws.Range("E46").Value = "  -3.51%  "

$ws.Range("E47").Value = "  -3.13%  "

$ws.Range("E48").Value = "  -2.59%  "

$ws.Range("E49").Value = "  -4.44%  "

$ws.Range("D50").Value = "'8.13"
$ws.Range("E50").Value = "  -5.46%  "

$ws.Range("D51").Value = "'132.21"
$ws.Range("E51").Value = "  -5.70%  "
